$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.646.98"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = "3.186.65"
$ws.Range("E3").Value = "  -3.93%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.34"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.60"
$ws.Range("E6").Value = "  -6.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").Value = "  -5.74%  "
$ws.Range("D9").Value = "3.196.28"
$ws.Range("E9").Value = "  -3.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.121"
$ws.Range("E10").Value = "  -3.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.86"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.394"
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").Value = "3.742.66"
$ws.Range("E13").Value = "  -3.92%  "
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "64.663.64"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.62"
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("D18").Value = "3.189.22"
$ws.Range("E18").Value = "  -2.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "421.54"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.04"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.35"
$ws.Range("E21").Value = "  -3.46%  "
$ws.Range("E22").Value = "  -2.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.45"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.206"
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.500"
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000106"
$ws.Range("E28").Value = "  -6.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.85"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.996"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("E31").Value = "  -3.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.85"
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.09"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.40"
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.13"
$ws.Range("E36").Value = "  -3.53%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.39"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.37"
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.71"
$ws.Range("E39").Value = "  -4.78%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.714.20"
$ws.Range("E40").Value = "  -5.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.27"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "24.43"
$ws.Range("E42").Value = "  -6.93%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.721"
$ws.Range("E43").Value = "  -4.86%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.16"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0626"
$ws.Range("E45").Value = "  -5.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.63"
$ws.Range("E46").Value = "  -4.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0264"
$ws.Range("E47").Value = "  -3.23%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.65"
$ws.Range("E48").Value = "  -6.52%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "293.37"
$ws.Range("E49").Value = "  -6.49%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.02"
$ws.Range("E50").Value = "  -11.82%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0993"
$ws.Range("E51").Value = "  -5.92%  "
